$wb = $excel.ActiveWorkbook

# "commands to remember" sheet: add a new row documenting the custom
# settings.xml maven build command, and make this the active sheet/selection.
$ws = $wb.Worksheets.Item("commands to remember")

$ws.Range("A4").Value = "mvn clean install "
$ws.Range("B4").Value = 'mvn clean install -s"company_resources_management_application_settings.xml"'
$ws.Range("C4").Value = "Cleans the older target file and reinstalls and builds the application from beginning when the settings.xml file is present inside the appication folder "

$ws.Activate()
$ws.Range("B12").Select() | Out-Null
